$d = $word.ActiveDocument

# --- Edit 1: Append red-colored "(This is a change - Version for main branch)"
#     to the end of the first paragraph, after two extra trailing spaces. ---
$p1 = $d.Paragraphs(1).Range
$end1 = $p1.End - 1                      # before the paragraph mark
$r = $d.Range($end1, $end1)
$r.InsertAfter("  ")

$r2start = $end1 + 2
$r2 = $d.Range($r2start, $r2start)
$r2.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$r2.Font.Color = 255                     ; # wdColorRed

$r3start = $r2.End
$r3 = $d.Range($r3start, $r3start)
$r3.InsertAfter("rsion for main branch")
$r3.Font.Color = 255

$r4start = $r3.End
$r4 = $d.Range($r4start, $r4start)
$r4.InsertAfter(")")
$r4.Font.Color = 255

# --- Edit 2: Remove the final paragraph ("...ank God almighty, we are free at last.") ---
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex).Range
$lastPara.Delete()

# --- Edit 3: Drop styles that are no longer used now that the removed
#     paragraph (and its NormalWeb-based siblings) are gone. Deleting from
#     the highest index down avoids re-resolving shifted indices. ---
$styleIndicesToDelete = @(18, 17, 16, 15, 14, 13, 12, 11, 10, 3, 2)
foreach ($idx in $styleIndicesToDelete) {
    $d.Styles($idx).Delete()
}
